# Update "Estado de Cuenta" worksheet: shift the period (Periodo Mora) table from
# Jul2025-Jan2025 (descending, 7 periods) to Mar2025-Aug2025 (ascending, 6 periods),
# and refresh the summary totals accordingly (adds "parte 1" of the new account statement data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Replace the period labels in the data table (rows 16-20 keep their old formatting,
#    row 21 will become the new last row after row 22 is removed below).
$ws.Range("E16").Value = "2503"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2506"
$ws.Range("E20").Value = "2507"

# 2. Row 21 becomes the last row of the table once row 22 is deleted, so give it the
#    same ("closing") border formatting that row 22 (the old last row) currently has.
$src = $ws.Range("B22:J22")
$dst = $ws.Range("B21:J21")
$src.Copy() | Out-Null
$dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = 0

# 3. Set the new (6th) period label for what is now the last data row.
$ws.Range("E21").Value = "2508"

# 4. Remove the now-obsolete oldest period row (period 2501); rows below move up
#    automatically (signature rows 27/28 -> 26/27).
$ws.Rows.Item(22).Delete() | Out-Null

# 5. Refresh the summary figures at the top: total "Valor Mora" and "Cant. Periodos".
$ws.Range("E11").Value = 341640
$ws.Range("F13").Value = 6
